# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by copying "2021-Q4" (same 8-col
#    fund-holding layout), inserted right after it and before "总计".
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newSheet = $wb.Worksheets.Item($q4.Index + 1)
$newSheet.Name = "2022-Q1"

# Force text storage (keep leading zeros / decimal formatting) for the
# numeric-looking string columns, matching the data's original inlineStr text.
$newSheet.Range("B2:G7").NumberFormat = "@"

$fundRows = @(
    @("012930", "中庚价值先锋股票", "54.59", "94.46", "4.40", "2.4020", 7),
    @("004450", "嘉实前沿科技沪港深股票", "21.72", "88.32", "5.38", "1.1685", 6),
    @("161914", "万家创业板2年定期开放混合A", "15.74", "95.20", "5.21", "0.8201", 10),
    @("161915", "万家创业板2年定期开放混合C", "2.36", "95.20", "5.21", "0.1230", 10),
    @("005104", "富荣福康混合A", "0.08", "87.88", "2.98", "0.0024", 10),
    @("005105", "富荣福康混合C", "0.04", "87.88", "2.98", "0.0012", 10)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ------------------------------------------------------------------
# 2. Prepend a 2022-Q1 summary row to the "总计" sheet, shifting the
#    existing rows down by one and renumbering the column-A index.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# re-apply the per-column formatting of the data rows (not the header row)
# to the freshly inserted row, since Insert() otherwise pulls it from row 1
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 6
$total.Cells.Item(2, 4).Value = 4.52

for ($r = 3; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# restore the originally-active tab (copying/editing sheets shifts focus)
$wb.Worksheets.Item("2020-Q4").Activate()
